# Apply the data refresh changes to both the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # F2: 想去人数 8860 -> 8893
    $ws.Range("F2").Value = 8893

    # F3: 想去人数 200 -> 201
    $ws.Range("F3").Value = 201

    # F4: 想去人数 433 -> 438
    $ws.Range("F4").Value = 438

    # I4: Cover image URL updated
    $ws.Range("I4").Value = "//i1.hdslb.com/bfs/openplatform/202402/FSJIeLNT1707130460798.jpeg"

    # F5: 想去人数 446 -> 448
    $ws.Range("F5").Value = 448
}
